$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 143; existing rows 143-151 shift down to 146-154.
$ws.Rows("143:145").Insert()

# New row 143: Especial, volumen 300, precios 23000, Curicó, 1533
$ws.Cells.Item(143, 1).Value = 4
$ws.Cells.Item(143, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(143, 3).Value = "Los Lagos"
$ws.Cells.Item(143, 4).Value = 44509
$ws.Cells.Item(143, 5).Value = 10
$ws.Cells.Item(143, 6).Value = "Fruta"
$ws.Cells.Item(143, 7).Value = 100101
$ws.Cells.Item(143, 8).Value = "Berries"
$ws.Cells.Item(143, 9).Value = 100101007
$ws.Cells.Item(143, 10).Value = "Kiwi"
$ws.Cells.Item(143, 11).Value = "Hayward"
$ws.Cells.Item(143, 12).Value = "Especial"
$ws.Cells.Item(143, 13).Value = 300
$ws.Cells.Item(143, 14).Value = 23000
$ws.Cells.Item(143, 15).Value = 23000
$ws.Cells.Item(143, 16).Value = 23000
$ws.Cells.Item(143, 17).Value = "$/caja 15 kilos"
$ws.Cells.Item(143, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(143, 19).Value = 1533
$ws.Cells.Item(143, 20).Value = 15

# New row 144: Primera, volumen 300, precios 17000, Curicó, 1133
$ws.Cells.Item(144, 1).Value = 4
$ws.Cells.Item(144, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(144, 3).Value = "Los Lagos"
$ws.Cells.Item(144, 4).Value = 44509
$ws.Cells.Item(144, 5).Value = 10
$ws.Cells.Item(144, 6).Value = "Fruta"
$ws.Cells.Item(144, 7).Value = 100101
$ws.Cells.Item(144, 8).Value = "Berries"
$ws.Cells.Item(144, 9).Value = 100101007
$ws.Cells.Item(144, 10).Value = "Kiwi"
$ws.Cells.Item(144, 11).Value = "Hayward"
$ws.Cells.Item(144, 12).Value = "Primera"
$ws.Cells.Item(144, 13).Value = 300
$ws.Cells.Item(144, 14).Value = 17000
$ws.Cells.Item(144, 15).Value = 17000
$ws.Cells.Item(144, 16).Value = 17000
$ws.Cells.Item(144, 17).Value = "$/caja 15 kilos"
$ws.Cells.Item(144, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(144, 19).Value = 1133
$ws.Cells.Item(144, 20).Value = 15

# New row 145: Segunda, volumen 300, precios 16000, Curicó, 1067
$ws.Cells.Item(145, 1).Value = 4
$ws.Cells.Item(145, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(145, 3).Value = "Los Lagos"
$ws.Cells.Item(145, 4).Value = 44509
$ws.Cells.Item(145, 5).Value = 10
$ws.Cells.Item(145, 6).Value = "Fruta"
$ws.Cells.Item(145, 7).Value = 100101
$ws.Cells.Item(145, 8).Value = "Berries"
$ws.Cells.Item(145, 9).Value = 100101007
$ws.Cells.Item(145, 10).Value = "Kiwi"
$ws.Cells.Item(145, 11).Value = "Hayward"
$ws.Cells.Item(145, 12).Value = "Segunda"
$ws.Cells.Item(145, 13).Value = 300
$ws.Cells.Item(145, 14).Value = 16000
$ws.Cells.Item(145, 15).Value = 16000
$ws.Cells.Item(145, 16).Value = 16000
$ws.Cells.Item(145, 17).Value = "$/caja 15 kilos"
$ws.Cells.Item(145, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(145, 19).Value = 1067
$ws.Cells.Item(145, 20).Value = 15
